$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / account holder info ---
$ws.Range("C2").Value = "Hartmut"
# Card number is a long digit string that must stay text (like the original cell);
# leading apostrophe forces text entry instead of numeric conversion.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 25.09.2023"

# --- Row 6 (existing transaction, update dates/desc/amount) ---
$ws.Range("B6").Value = "27.09."
$ws.Range("C6").Value = "28.09."
$ws.Range("D6").Value = "AMAZON.DE MKTPLC EU DWLNBD"
$ws.Range("E6").Value = "188,34-"

# --- Row 7 (existing transaction, dates/amount change, description same) ---
$ws.Range("B7").Value = "29.09."
$ws.Range("C7").Value = "30.09."
$ws.Range("E7").Value = "24,60-"

# --- Row 8 (existing transaction, update dates/desc/amount) ---
$ws.Range("B8").Value = "02.10."
$ws.Range("C8").Value = "03.10."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 43408064"
$ws.Range("E8").Value = "39,26-"

# --- Row 9 (was blank, now a new transaction row) ---
# Copy formatting from row 8 amount cell so E9 matches style s=17 used by other amount cells
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B9").Value = "03.10."
$ws.Range("C9").Value = "04.10."
$ws.Range("D9").Value = "BURGER KING Gifhorn"
$ws.Range("E9").Value = "31,83-"

# --- Row 10 (was blank, now a new transaction row) ---
$ws.Range("E8").Copy()
$ws.Range("E10").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B10").Value = "05.10."
$ws.Range("C10").Value = "06.10."
$ws.Range("D10").Value = "BEITRAG Allianz SE K-90949889"
$ws.Range("E10").Value = "55,67-"

# --- Row 11 (was blank, now a new transaction row) ---
$ws.Range("E8").Copy()
$ws.Range("E11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B11").Value = "07.10."
$ws.Range("C11").Value = "08.10."
$ws.Range("D11").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 39557599"
$ws.Range("E11").Value = "86,92-"

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 09.10.2023"
$ws.Range("E12").Value = "426,62-"

# --- Next billing date note ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 18.10.2023"
